# Update the "Protocol: Context Driven Interaction ..." bullet with the
# expanded wording, then remove the now-redundant blank paragraph that
# used to follow it.

$d = $word.ActiveDocument

$old = "Protocol: Context Driven Interaction REST P2P (SIDs CDI Dialogs: runat peer resolution semantics)."
$new = "Protocol: Context Driven Interaction REST P2P (SIDs CDI Dialogs: runat peer resolution addressable / browseable messages interactions embedded session semantics: events sourcing / history terms resolution)."

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find target paragraph text to replace."
}

# Locate the paragraph that now holds the updated text, then delete the
# following empty paragraph (the one with the pBdr/shd/ind formatting and
# no run text) so the two paragraphs collapse into one, matching the diff.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Protocol: Context Driven Interaction REST P2P*") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text.Trim() -eq "") {
            $next.Range.Delete()
        }
        break
    }
}
